$wb = $excel.ActiveWorkbook

$newPlants = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

# --- CRtPaL-profits sheet: new plant types get a "change in capacity" value of 2 ---
$wsProfits = $wb.Worksheets.Item("CRtPaL-profits")
$row = 19
foreach ($plant in $newPlants) {
    $wsProfits.Range("A$row").Value = $plant
    $wsProfits.Range("B$row").Value = 2
    $wsProfits.Range("B$row").NumberFormat = "0"
    $row = $row + 1
}
$wsProfits.Range("A19:A24").Select()

# --- CRtPaL-losses sheet: new plant types get a "change in capacity" value of 1 ---
$wsLosses = $wb.Worksheets.Item("CRtPaL-losses")
$row = 19
foreach ($plant in $newPlants) {
    $wsLosses.Range("A$row").Value = $plant
    $wsLosses.Range("B$row").Value = 1
    $wsLosses.Range("B$row").NumberFormat = "0"
    $row = $row + 1
}
$wsLosses.Range("A19:A24").Select()
